$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 23:22"

# --- Row 7: Francia ---
$ws.Range("B7").Value = 166420
$ws.Range("C7").Value = 509
$ws.Range("E7").Value = 94105

# --- Row 14: Brasil ---
$ws.Range("B14").Value = 78162
$ws.Range("C14").Value = 5263
$ws.Range("D14").Value = 34132
$ws.Range("E14").Value = 38564
$ws.Range("G14").Value = 403
$ws.Range("H14").Value = 5466

# --- Rows 18/19: Peru overtakes India, so Peru now sits above India ---
# Row 18 becomes Peru with fresh data, row 19 becomes India with Peru's old slot data
$ws.Range("A18").Value = "Peru"
$ws.Range("B18").Value = 33931
$ws.Range("C18").Value = 2741
$ws.Range("D18").Value = 10037
$ws.Range("E18").Value = 22951
$ws.Range("F18").Value = 623
$ws.Range("G18").Value = 89
$ws.Range("H18").Value = 943

$ws.Range("A19").Value = "India"
$ws.Range("B19").Value = 33062
$ws.Range("C19").Value = 1738
$ws.Range("D19").Value = 8437
$ws.Range("E19").Value = 23546
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 71
$ws.Range("H19").Value = 1079
